# Update TPM-derived NATMI values for Tnfsf10-Tnfrsf10b LR pair
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.31937166666667
$ws.Range("H2").Value = 51.958115
$ws.Range("I2").Value = 0.9799257492601914
$ws.Range("J2").Value = 0.9799257492601915
$ws.Range("M2").Value = 5.857771333333333
$ws.Range("N2").Value = 17.573314
$ws.Range("O2").Value = 0.7247832978729415
$ws.Range("P2").Value = 0.7247832978729416
$ws.Range("Q2").Value = 101.4529188603455
$ws.Range("R2").Value = 913.0762697431099
$ws.Range("S2").Value = 0.7102338162194146
$ws.Range("T2").Value = 0.7102338162194148

$ws.Range("G3").Value = 17.31937166666667
$ws.Range("H3").Value = 51.958115
$ws.Range("I3").Value = 0.9799257492601914
$ws.Range("J3").Value = 0.9799257492601915
$ws.Range("O3").Value = 0.1568878485835839
$ws.Range("P3").Value = 0.1568878485835839
$ws.Range("Q3").Value = 21.96067461713888
$ws.Range("R3").Value = 197.64607155425
$ws.Range("S3").Value = 0.1537384425730879
$ws.Range("T3").Value = 0.1537384425730879

$ws.Range("G4").Value = 17.31937166666667
$ws.Range("H4").Value = 51.958115
$ws.Range("I4").Value = 0.9799257492601914
$ws.Range("J4").Value = 0.9799257492601915
$ws.Range("O4").Value = 0.1183288535434745
$ws.Range("P4").Value = 0.1183288535434745
$ws.Range("Q4").Value = 16.56330604280611
$ws.Range("R4").Value = 149.069754385255
$ws.Range("S4").Value = 0.1159534904676887
$ws.Range("T4").Value = 0.1159534904676887

$ws.Range("I5").Value = 0.009666424266436919
$ws.Range("J5").Value = 0.009666424266436919
$ws.Range("M5").Value = 5.857771333333333
$ws.Range("N5").Value = 17.573314
$ws.Range("O5").Value = 0.7247832978729415
$ws.Range("P5").Value = 0.7247832978729416
$ws.Range("Q5").Value = 1.000776801214667
$ws.Range("R5").Value = 9.006991210932
$ws.Range("S5").Value = 0.00700606285846718
$ws.Range("T5").Value = 0.007006062858467181

$ws.Range("I6").Value = 0.009666424266436919
$ws.Range("J6").Value = 0.009666424266436919
$ws.Range("O6").Value = 0.1568878485835839
$ws.Range("P6").Value = 0.1568878485835839
$ws.Range("S6").Value = 0.001516544506657437
$ws.Range("T6").Value = 0.001516544506657437

$ws.Range("I7").Value = 0.009666424266436919
$ws.Range("J7").Value = 0.009666424266436919
$ws.Range("O7").Value = 0.1183288535434745
$ws.Range("P7").Value = 0.1183288535434745
$ws.Range("S7").Value = 0.001143816901312302
$ws.Range("T7").Value = 0.001143816901312302

$ws.Range("I8").Value = 0.01040782647337163
$ws.Range("J8").Value = 0.01040782647337163
$ws.Range("M8").Value = 5.857771333333333
$ws.Range("N8").Value = 17.573314
$ws.Range("O8").Value = 0.7247832978729415
$ws.Range("P8").Value = 0.7247832978729416
$ws.Range("Q8").Value = 1.077535084176222
$ws.Range("R8").Value = 9.697815757586001
$ws.Range("S8").Value = 0.007543418795059594
$ws.Range("T8").Value = 0.007543418795059597

$ws.Range("I9").Value = 0.01040782647337163
$ws.Range("J9").Value = 0.01040782647337163
$ws.Range("O9").Value = 0.1568878485835839
$ws.Range("P9").Value = 0.1568878485835839
$ws.Range("S9").Value = 0.001632861503838544
$ws.Range("T9").Value = 0.001632861503838544

$ws.Range("I10").Value = 0.01040782647337163
$ws.Range("J10").Value = 0.01040782647337163
$ws.Range("O10").Value = 0.1183288535434745
$ws.Range("P10").Value = 0.1183288535434745
$ws.Range("S10").Value = 0.001231546174473488
$ws.Range("T10").Value = 0.001231546174473488
